# Refresh the cryptocurrency Price (D) and Volume(1h) (E) columns on Sheet1
# with the latest scraped values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.413.04'
$ws.Range("E2").Value = '  +4.27%  '

$ws.Range("D3").Value = '2.429.23'
$ws.Range("E3").Value = '  +5.60%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.51'
$ws.Range("E5").Value = '  +2.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.72'
$ws.Range("E6").Value = '  +7.31%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("E8").Value = '  +2.20%  '

$ws.Range("D9").Value = '2.427.90'
$ws.Range("E9").Value = '  +5.62%  '

$ws.Range("E10").Value = '  +3.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.75'
$ws.Range("E11").Value = '  +4.21%  '

$ws.Range("E12").Value = '  +0.47%  '

$ws.Range("E13").Value = '  +5.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.29'
$ws.Range("E14").Value = '  +13.01%  '

$ws.Range("D15").Value = '2.862.08'
$ws.Range("E15").Value = '  +5.63%  '

$ws.Range("D16").Value = '62.300.43'
$ws.Range("E16").Value = '  +4.11%  '

$ws.Range("E17").Value = '  +7.07%  '

$ws.Range("D18").Value = '2.425.84'
$ws.Range("E18").Value = '  +5.73%  '

$ws.Range("E19").Value = '  +7.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.23'
$ws.Range("E20").Value = '  +11.31%  '

$ws.Range("E21").Value = '  +3.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.82'
$ws.Range("E22").Value = '  +4.18%  '

$ws.Range("E23").Value = '  +0.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.14'
$ws.Range("E24").Value = '  +2.29%  '

$ws.Range("E25").Value = '  +1.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.53'
$ws.Range("E27").Value = '  +14.39%  '

$ws.Range("E28").Value = '  +5.71%  '

$ws.Range("E29").Value = '  +14.46%  '

$ws.Range("E30").Value = '  +5.57%  '

$ws.Range("D31").Value = '0.0₃0785'
$ws.Range("E31").Value = '  +8.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.42'
$ws.Range("E32").Value = '  +10.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '171.93'
$ws.Range("E33").Value = '  +0.91%  '

$ws.Range("E34").Value = '  +6.66%  '

$ws.Range("E35").Value = '  +5.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '381.60'
$ws.Range("E36").Value = '  +20.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.55'
$ws.Range("E37").Value = '  +5.19%  '

$ws.Range("E38").Value = '  +11.91%  '

$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("E40").Value = '  -0.08%  '

$ws.Range("E41").Value = '  +12.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.18'
$ws.Range("E42").Value = '  +3.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '145.04'
$ws.Range("E43").Value = '  +6.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.66'
$ws.Range("E44").Value = '  +7.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.69'
$ws.Range("E45").Value = '  +10.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.591'
$ws.Range("E46").Value = '  +5.17%  '

$ws.Range("E48").Value = '  +6.13%  '

$ws.Range("E49").Value = '  +5.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.89'
$ws.Range("E50").Value = '  +7.04%  '

$ws.Range("D51").Value = '0.0₆0216'
$ws.Range("E51").Value = '  -1.98%  '
